$d = $word.ActiveDocument

# Locate the target paragraph ("About Me" section, the sentence about
# education / spare time / mentoring) by finding the unique anchor text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*I value education*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

# Clear the paragraph's existing text, keeping the trailing paragraph mark.
$full = $target.Range.Duplicate
$full.MoveEnd(1, -1) | Out-Null
$full.Delete()

# Remove the old _GoBack bookmark (it currently sits at the end of the
# paragraph); we'll re-create it in its new location below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$start = $target.Range.Start

# Insert the full replacement text as a single run.
$ip = $d.Range($start, $start)
$ip.InsertBefore("I value education for the pursuit of knowledge and development of my skills. In my spare time I enjoy leading an advocacy team for a student group on campus. ")

# Re-create the _GoBack bookmark right after "In my spare time" - this
# also splits the text into two runs at that point.
$bmOffset = $start + "I value education for the pursuit of knowledge and development of my skills. In my spare time".Length
$bmRange = $d.Range($bmOffset, $bmOffset)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Split the remaining text into separate runs matching the target
# structure, using temporary bookmarks as run-boundary markers (adding
# then immediately deleting a bookmark splits the underlying run without
# leaving the bookmark behind).
$splitPoints = @(
    " I enjoy ".Length,
    "leading an".Length,
    " advocacy ".Length,
    "team".Length
)

$offset = $bmOffset
foreach ($len in $splitPoints) {
    $offset += $len
    $splitRange = $d.Range($offset, $offset)
    $d.Bookmarks.Add("TempSplit", $splitRange)
    $d.Bookmarks("TempSplit").Delete()
}

Write-Output ("Final paragraph text: [" + $target.Range.Text + "]")
